$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add row 6: new status update "25/08/2016" with progress through Compose Mail
$ws.Range("A6").Value = "25/08/2016"
$ws.Range("B6").Value = "done"
$ws.Range("C6").Value = "done"
$ws.Range("D6").Value = "done"
$ws.Range("E6").Value = "done"
$ws.Range("F6").Value = "done"
$ws.Range("G6").Value = "done"
$ws.Range("H6").Value = "done"
$ws.Range("I6").Value = "working in progress"

$ws.Range("I6").Select()
